$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '48.213.95'
$ws.Cells.Item(2, 5).Value = '  +2.28%  '

$ws.Cells.Item(3, 4).Value = '2.526.77'
$ws.Cells.Item(3, 5).Value = '  +1.52%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '323.90'
$ws.Range('D5').Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.59%  '

$ws.Cells.Item(6, 5).Value = '  +0.65%  '

$ws.Cells.Item(7, 5).Value = '  +0.91%  '

$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.557'
$ws.Range('D9').Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +4.18%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '40.99'
$ws.Range('D10').Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +5.68%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '20.48'
$ws.Range('D11').Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +11.88%  '

$ws.Cells.Item(13, 5).Value = '  +1.38%  '

$ws.Cells.Item(14, 5).Value = '  +1.78%  '

$ws.Cells.Item(15, 4).Value = '2.922.34'
$ws.Cells.Item(15, 5).Value = '  +1.60%  '

$ws.Cells.Item(16, 4).Value = '2.530.82'
$ws.Cells.Item(16, 5).Value = '  +1.57%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.860'
$ws.Range('D17').Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +1.17%  '

$ws.Cells.Item(18, 4).Value = '48.057.95'
$ws.Cells.Item(18, 5).Value = '  +2.14%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '13.31'
$ws.Range('D19').Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +4.83%  '

$ws.Cells.Item(20, 5).Value = '  +0.49%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0950'
$ws.Cells.Item(21, 5).Value = '  +1.41%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.71'
$ws.Range('D22').Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.84%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '72.19'
$ws.Range('D23').Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +2.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '269.16'
$ws.Range('D24').Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +9.24%  '

$ws.Cells.Item(25, 5).Value = '  +0.43%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '26.23'
$ws.Range('D26').Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +1.56%  '

$ws.Cells.Item(27, 5).Value = '  -0.26%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.18'
$ws.Range('D28').Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.61%  '

$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.145'
$ws.Range('D29').Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +3.24%  '

$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '2.21'
$ws.Range('D30').Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -3.47%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '35.84'
$ws.Range('D31').Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +2.30%  '

$ws.Cells.Item(32, 5).Value = '  -0.64%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '20.00'
$ws.Range('D33').Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +0.08%  '

$ws.Cells.Item(34, 5).Value = '  +0.26%  '

$ws.Cells.Item(35, 5).Value = '  +0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0795'
$ws.Range('D36').Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.50%  '

$ws.Cells.Item(37, 5).Value = '  +1.66%  '

$ws.Cells.Item(38, 5).Value = '  +1.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.01'
$ws.Range('D39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.113'
$ws.Range('D40').Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.36%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '22.16'
$ws.Range('D41').Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +4.47%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '120.22'
$ws.Range('D42').Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -0.61%  '

$ws.Cells.Item(43, 5).Value = '  -1.73%  '

$ws.Cells.Item(44, 5).Value = '  +1.85%  '

$ws.Cells.Item(45, 4).Value = '2.021.90'
$ws.Cells.Item(45, 5).Value = '  +1.64%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.17'
$ws.Range('D46').Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +4.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.91'
$ws.Range('D47').Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +6.78%  '

$ws.Cells.Item(48, 5).Value = '  +0.33%  '

$ws.Cells.Item(49, 5).Value = '  +0.77%  '

$ws.Cells.Item(50, 5).Value = '  +2.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '79.87'
$ws.Range('D51').Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +3.18%  '
